# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff file id
# (44ea2c0c-9e7f-4991-920f-82ee823ce333 replacing
#  f3621c04-1fb4-4a02-9775-7e3ee7e7defb) and refreshed handoff timestamps,
# and clears the (not-yet-produced) handback info for zh-cn / de-de.

$wb = $excel.ActiveWorkbook

$oldId = "f3621c04-1fb4-4a02-9775-7e3ee7e7defb"
$newId = "44ea2c0c-9e7f-4991-920f-82ee823ce333"
$oldHash = "4453b1d97e2c9b9d5257939f998595323278fe14"
$newHash = "3af080a4d18a6a288d2c16290b2e5bc9e0b2a98e"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newId.md"

$ov.Hyperlinks.Delete()
$ov.Range("B2").Value = "e2e\$newId.md"
$ovLink = $ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d20118c3d65992a88bca8e5eb3d5335085188920/e2e/$newId.md", "", "", "e2e\$newId.md")

$ov.Range("G2").Value = "2016-09-04 01:05:35"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()
$zh.Range("A2").Value = "$newId.md"
$zhLink = $zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d20118c3d65992a88bca8e5eb3d5335085188920/e2e/$newId.md", "", "", "$newId.md")

$zh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-04 01:05:31"
$zh.Range("I2").Value = "'"
$zh.Range("J2").Value = "'"
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()
$de.Range("A2").Value = "$newId.md"
$deLink = $de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d20118c3d65992a88bca8e5eb3d5335085188920/e2e/$newId.md", "", "", "$newId.md")

$de.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$de.Range("H2").Value = "2016-09-04 01:05:35"
$de.Range("I2").Value = "'"
$de.Range("J2").Value = "'"
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426

Write-Host "Done"
